$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the "keyword_order" (portal search ranking) row - this feature is
# being dropped from the page-info table.
$ws.Rows.Item(3).Delete()

# Insert a new blank row where the "health" row used to be (row 4), pushing
# passwd/fav_link/today_contents back down to their original row numbers.
$ws.Rows.Item(4).Insert()

# The freshly inserted row inherits the wrap-text style used in columns J:L on
# neighbouring rows; the new row has no SQL list/update/insert formats, so
# clear that stale formatting entirely (column M keeps its style below).
$ws.Range("J4:L4").Clear()

# Fill in the new "create_diary" row (order chosen to mirror the shared-string
# insertion order used when this edit was originally authored in Excel).
$ws.Range("C4").Value = "원석의 일기"
$ws.Range("E4").Value = "/create_diary.neo"
$ws.Range("A4").Value = "create_diary"
$ws.Range("B4").Value = "CreateDiaryWebApp"
$ws.Range("H4").Value = "개인 일기 제목 자동 생성 페이지"
$ws.Range("D4").Value = "create_diary.html"
$ws.Range("F4").Value = "kwo"
$ws.Range("G4").Value = "keyword_order"
$ws.Range("I4").Value = "검색어"
$ws.Range("M4").Value = "private"

# Update the active selection to match the freeze-pane view used after the edit.
[void]$ws.Activate()
[void]$ws.Range("D4").Select()
